# Auto-generated edit script applying scraped market-data updates
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H18").Value = 1433
$ws.Range("I18").Value = 1433
$ws.Range("K18").Value = 1433
$ws.Range("M18").Value = -1149

$ws.Range("H34").Value = 4992
$ws.Range("I34").Value = 4992
$ws.Range("K34").Value = 4992
$ws.Range("M34").Value = -4789

$ws.Range("H36").Value = 4992
$ws.Range("I36").Value = 4992
$ws.Range("K36").Value = 4992
$ws.Range("M36").Value = -4277

$ws.Range("H40").Value = 1998.8667
$ws.Range("I40").Value = 1808.6
$ws.Range("K40").Value = 1808.6
$ws.Range("M40").Value = -1633.6

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents() | Out-Null

$ws.Range("H51").Value = 6929.7
$ws.Range("I51").Value = 6600
$ws.Range("J51").Value = 7071
$ws.Range("K51").Value = 6600
$ws.Range("L51").Value = 7071
$ws.Range("M51").Value = -6116
$ws.Range("N51").Value = -8039

$ws.Range("H92").Value = 348.57144
$ws.Range("I92").Value = 89.666664
$ws.Range("J92").Value = 542.75
$ws.Range("K92").Value = 89.666664
$ws.Range("L92").Value = 542.75
$ws.Range("M92").Value = 1158.333336
$ws.Range("N92").Value = -3038.75

$ws.Range("H98").Value = 651.2941
$ws.Range("I98").Value = 643.25
$ws.Range("K98").Value = 643.25
$ws.Range("M98").Value = 854.75

$ws.Range("H106").Value = 9999.833000000001
$ws.Range("I106").Value = 8249.75
$ws.Range("K106").Value = 8249.75
$ws.Range("M106").Value = -7618.75

$ws.Range("H122").Value = 651.2941
$ws.Range("I122").Value = 643.25
$ws.Range("K122").Value = 1929.75
$ws.Range("M122").Value = 520.25

$ws.Range("H127").Value = 2525.25
$ws.Range("I127").Value = 2728.2222
$ws.Range("J127").Value = 1916.3334
$ws.Range("K127").Value = 8184.6666
$ws.Range("L127").Value = 5749.0002
$ws.Range("M127").Value = -3224.6666
$ws.Range("N127").Value = -15669.0002

$ws.Range("H132").Value = 3480.72
$ws.Range("I132").Value = 2281.05
$ws.Range("K132").Value = 6843.150000000001
$ws.Range("M132").Value = -4313.150000000001

$ws.Range("H138").Value = 2543.3333
$ws.Range("I138").Value = 2753.111
$ws.Range("J138").Value = 2228.6667
$ws.Range("K138").Value = 8259.332999999999
$ws.Range("L138").Value = 6686.000100000001
$ws.Range("M138").Value = -3119.332999999999
$ws.Range("N138").Value = -16966.0001

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H45").Value = 7950
$ws.Range("I45").Value = 7950
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 7950
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -7573
$ws.Range("N45").ClearContents() | Out-Null

$ws.Range("H61").Value = 2381.077
$ws.Range("I61").Value = 2276.32
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2276.32
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2064.32
$ws.Range("N61").Value = -5424

$ws.Range("H88").Value = 799.5714
$ws.Range("J88").Value = 868.25
$ws.Range("L88").Value = 868.25
$ws.Range("N88").Value = -1680.25

$ws.Range("H91").Value = 799.5714
$ws.Range("J91").Value = 868.25
$ws.Range("L91").Value = 868.25
$ws.Range("N91").Value = -3676.25

$ws.Range("H97").Value = 1208.875
$ws.Range("I97").Value = 1208.875
$ws.Range("K97").Value = 1208.875
$ws.Range("M97").Value = -712.875

$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080

$ws.Range("H132").Value = 310
$ws.Range("I132").Value = 310
$ws.Range("K132").Value = 930
$ws.Range("M132").Value = 1600

$ws.Range("H136").Value = 2381.077
$ws.Range("I136").Value = 2276.32
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6828.960000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4278.960000000001
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H22").Value = 313.85715
$ws.Range("I22").Value = 316.16666
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 316.16666
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -143.16666
$ws.Range("N22").Value = -646

$ws.Range("H86").Value = 1622.6111
$ws.Range("I86").Value = 1513.8
$ws.Range("J86").Value = 2166.6667
$ws.Range("K86").Value = 1513.8
$ws.Range("L86").Value = 2166.6667
$ws.Range("M86").Value = -390.8
$ws.Range("N86").Value = -4412.6667

$ws.Range("H89").Value = 1622.6111
$ws.Range("I89").Value = 1513.8
$ws.Range("J89").Value = 2166.6667
$ws.Range("K89").Value = 7569
$ws.Range("L89").Value = 10833.3335
$ws.Range("M89").Value = -1953
$ws.Range("N89").Value = -22065.3335

$ws.Range("H94").Value = 663.92
$ws.Range("I94").Value = 582.6087
$ws.Range("K94").Value = 582.6087
$ws.Range("M94").Value = -131.6087

$ws.Range("H105").Value = 4632.3335
$ws.Range("I105").Value = 4599
$ws.Range("J105").Value = 4699
$ws.Range("K105").Value = 4599
$ws.Range("L105").Value = 4699
$ws.Range("M105").Value = -2852
$ws.Range("N105").Value = -8193

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H20").Value = 44869
$ws.Range("J20").Value = 44869
$ws.Range("L20").Value = 44869
$ws.Range("N20").Value = -45341

$ws.Range("H30").Value = 44869
$ws.Range("J30").Value = 44869
$ws.Range("L30").Value = 44869
$ws.Range("N30").Value = -45051

$ws.Range("H128").Value = 44869
$ws.Range("J128").Value = 44869
$ws.Range("L128").Value = 44869
$ws.Range("N128").Value = -54829

$ws.Range("H132").Value = 1961.4286
$ws.Range("I132").Value = 2022.2354
$ws.Range("J132").Value = 1867.4546
$ws.Range("K132").Value = 6066.706200000001
$ws.Range("L132").Value = 5602.3638
$ws.Range("M132").Value = -3536.706200000001
$ws.Range("N132").Value = -10662.3638

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H3").Value = 8599.6
$ws.Range("I3").Value = 6499
$ws.Range("K3").Value = 19497
$ws.Range("M3").Value = -19385

$ws.Range("H11").Value = 26804682
$ws.Range("I11").Value = 34114980
$ws.Range("K11").Value = 102344940
$ws.Range("M11").Value = -102344800

$ws.Range("H121").Value = 10293.647
$ws.Range("I121").Value = 23295.6
$ws.Range("J121").Value = 4876.1665
$ws.Range("K121").Value = 69886.79999999999
$ws.Range("L121").Value = 14628.4995
$ws.Range("M121").Value = -68576.79999999999
$ws.Range("N121").Value = -17248.4995

$ws.Range("H122").Value = 685.6923
$ws.Range("J122").Value = 709.5
$ws.Range("L122").Value = 6385.5
$ws.Range("N122").Value = -11285.5

$ws.Range("H131").Value = 1002313.3
$ws.Range("I131").Value = 1050
$ws.Range("J131").Value = 1252629.1
$ws.Range("K131").Value = 3150
$ws.Range("L131").Value = 3757887.3
$ws.Range("M131").Value = 1890
$ws.Range("N131").Value = -3767967.3

$ws.Range("H132").Value = 2542.8
$ws.Range("J132").Value = 3847.25
$ws.Range("L132").Value = 34625.25
$ws.Range("N132").Value = -39685.25

$ws.Range("H141").Value = 13265.5
$ws.Range("I141").Value = 14176.333
$ws.Range("K141").Value = 42528.999
$ws.Range("M141").Value = -37348.999

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H113").Value = 1676
$ws.Range("I113").Value = 1345
$ws.Range("K113").Value = 1345
$ws.Range("M113").Value = 825

$ws.Range("H132").Value = 1592.9231
$ws.Range("I132").Value = 1431.5714
$ws.Range("K132").Value = 4294.7142
$ws.Range("M132").Value = -1764.7142

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H46").Value = 102696
$ws.Range("I46").Value = 126870
$ws.Range("J46").Value = 6000
$ws.Range("K46").Value = 126870
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -126682
$ws.Range("N46").Value = -6376

$ws.Range("H75").Value = 500173
$ws.Range("J75").Value = 500173
$ws.Range("L75").Value = 500173
$ws.Range("N75").Value = -502045

$ws.Range("H78").Value = 500173
$ws.Range("J78").Value = 500173
$ws.Range("L78").Value = 1500519
$ws.Range("N78").Value = -1509879

$ws.Range("H82").Value = 1301.3846
$ws.Range("I82").Value = 1375.2727
$ws.Range("K82").Value = 1375.2727
$ws.Range("M82").Value = -1014.2727

$ws.Range("H85").Value = 1301.3846
$ws.Range("I85").Value = 1375.2727
$ws.Range("K85").Value = 1375.2727
$ws.Range("M85").Value = -127.2727

$ws.Range("H93").Value = 1261.5
$ws.Range("I93").Value = 1198.1538
$ws.Range("K93").Value = 1198.1538
$ws.Range("M93").Value = 49.84619999999995

$ws.Range("H122").Value = 5747.357
$ws.Range("I122").Value = 5372.0415
$ws.Range("J122").Value = 7999.25
$ws.Range("K122").Value = 16116.1245
$ws.Range("L122").Value = 23997.75
$ws.Range("M122").Value = -13666.1245
$ws.Range("N122").Value = -28897.75

$ws.Range("H136").Value = 6000
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15450
$ws.Range("N136").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H18").Value = 13249.5
$ws.Range("I18").Value = 9749
$ws.Range("J18").Value = 16750
$ws.Range("K18").Value = 9749
$ws.Range("L18").Value = 16750
$ws.Range("M18").Value = -9576
$ws.Range("N18").Value = -17096

$ws.Range("H51").Value = 32500
$ws.Range("J51").Value = 55000
$ws.Range("L51").Value = 55000
$ws.Range("N51").Value = -56020

$ws.Range("H130").Value = 37500
$ws.Range("J130").Value = 37500
$ws.Range("L130").Value = 37500
$ws.Range("N130").Value = -47540

$ws.Range("H132").Value = 5209.048
$ws.Range("I132").Value = 4670.4707
$ws.Range("K132").Value = 14011.4121
$ws.Range("M132").Value = -11481.4121

